$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24
$lastRow = 23

# Set the new row's values. The date column stores plain text like
# "2025/12/03" (same as the rows above it) rather than a real date, so a
# leading apostrophe forces text interpretation instead of Excel's
# automatic date parsing.
$ws.Cells.Item($newRow, 1).Value = "'2025/12/03"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1322

# Copy the formatting from the row above (centered alignment style) onto
# the new row so it matches the existing data rows.
$ws.Range("A$lastRow:C$lastRow").Copy() | Out-Null
$ws.Range("A${newRow}:C${newRow}").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
